$d = $word.ActiveDocument

# Locate the "Luu y:" paragraph via Find (robust against index drift).
$rng = $d.Content
$found = $rng.Find.Execute("Lưu ý:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    Write-Host "ERROR: could not find Luu y paragraph"
} else {
    $noteParaIndex = $rng.Paragraphs(1).Index
    $notePara = $d.Paragraphs($noteParaIndex)
    $prevPara = $notePara.Previous()
    $lastPara = $notePara.Next()

    # Copy the tab stop from the "Luu y" paragraph onto the preceding
    # ("Neu lai, in ra so tien lai...") paragraph before it disappears.
    $prevPara.TabStops.Add(56.7)

    # Remove the "Luu y: Tram dau tien..." and "Tram cuoi cung..." paragraphs
    # (and their paragraph marks) entirely.
    $delStart = $notePara.Range.Start
    $delEnd = $lastPara.Range.End
    $d.Range($delStart, $delEnd).Delete()

    Write-Host "OK"
}
